$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlinks (row 6: agtan@gmail.com / agtan4321) up front ---
$ws.Range("C6").Hyperlinks.Delete() | Out-Null

# --- Row 2: Nadia Hertisa Isnaeni Putri ---
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "Nadia Hertisa Isnaeni Putri"
$ws.Range("C2").Value = "hertisanadia44@gmail.com "
$ws.Range("D2").Value = "Komplek Permata Kopo C-189"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0043171547"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 192010523
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = "smkn1ktp@01"
$ws.Range("J2").Value = "user"

# --- Row 3: Arianti Apriani Sagita ---
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "Arianti Apriani Sagita"
$ws.Range("C3").Value = "ariantiaprianisagita@gmail.com"
$ws.Range("D3").Value = "Kp. Pasanggrahan Rt 02 Rw 06 Kec. Pasirjambu"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0023620702"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 192010505
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = "smkn1ktp@01"
$ws.Range("J3").Value = "user"

# --- Row 4: Ajeng Nurfadillah ---
$ws.Range("A4").Value = 14
$ws.Range("B4").Value = "Ajeng Nurfadillah"
$ws.Range("C4").Value = "ajengnurfadilah@gmail.com"
$ws.Range("D4").Value = "Jln. Raya Sayuran Rt 08 Rw 07"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0034169559"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 192010501
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = "smkn1ktp@01"
$ws.Range("J4").Value = "user"

# --- Row 5: Tiara Situmorang (A5 stays blank) ---
$ws.Range("A5").ClearContents() | Out-Null
$ws.Range("B5").Value = "Tiara Situmorang"

$ws.Range("C5").Value = "tiara222324@gmail.com"

$ws.Range("D5").Value = "Kp. Pamoyanan Rt03/Rw05"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Font.Name = "Calibri"
$ws.Range("E5").Font.Size = 11
$ws.Range("E5").Font.Color = 0
$ws.Range("E5").Value = "0045802564"

$ws.Range("F5").Value = 192010533
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 7

$ws.Range("I5").Value = "smkn1ktp@01"

$ws.Range("J5").Font.Name = "Calibri"
$ws.Range("J5").Font.Size = 11
$ws.Range("J5").Font.Color = 0
$ws.Range("J5").Value = "user"

# --- Remove the old row 6 entirely (Agtan Dwiputra) ---
$ws.Rows("6:6").Delete() | Out-Null

# --- New hyperlinks for row 5 ---
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:tiara222324@gmail.com") | Out-Null
$ws.Range("C5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("I5"), "mailto:smkn1ktp@01") | Out-Null
$ws.Range("I5").Style = "Hyperlink"

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to J6 (first empty row below the data) ---
$ws.Range("J6").Select() | Out-Null
